# Rename the inline logo pictures that live in the document's headers and
# footers:
#   - the BTEC logo (.jpg) pictures in the headers: image1.jpg -> image2.jpg
#   - the Pearson Edexcel logo (.png) pictures in the footers: image2.png -> image1.png
#
# Only the picture's display Name changes; the underlying media file and
# its relationship are left untouched.

$d = $word.ActiveDocument
$section = $d.Sections.Item(1)

# --- Headers (BTec_Logo-Orange): image1.jpg -> image2.jpg ------------------
$headerPrimary = $section.Headers.Item(1).Range.InlineShapes.Item(1)
$headerPrimary.Name = "image2.jpg"

$headerFirstPage = $section.Headers.Item(2).Range.InlineShapes.Item(1)
$headerFirstPage.Name = "image2.jpg"

# --- Footers (PearsonLogo.png): image2.png -> image1.png -------------------
# Re-seating the shape through its own .Range before setting .Name avoids a
# stale-handle resolution error against these multi-paragraph footer stories.
$footerPrimary = $section.Footers.Item(1).Range.InlineShapes.Item(1)
$footerPrimary = $footerPrimary.Range.InlineShapes.Item(1)
$footerPrimary.Name = "image1.png"

$footerFirstPage = $section.Footers.Item(2).Range.InlineShapes.Item(1)
$footerFirstPage = $footerFirstPage.Range.InlineShapes.Item(1)
$footerFirstPage.Name = "image1.png"
